$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 184-218 with changed cell values ---
$ws.Cells.Item(184, 5).Value = 'PAPELON'  # E184: was 'LIOCDEP'
$ws.Cells.Item(184, 8).Value = 22  # H184: was 24
$ws.Cells.Item(184, 9).Value = 7  # I184: was -1
$ws.Cells.Item(184, 10).Value = 'F'  # J184: was 'M'
$ws.Cells.Item(184, 11).Value = 'MEDPF-1'  # K184: was None
$ws.Cells.Item(184, 16).Value = 'NA l inferred'  # P184: was None
$ws.Cells.Item(185, 5).Value = 'PAPELON'  # E185: was 'LIOCDEP'
$ws.Cells.Item(185, 8).Value = 22  # H185: was 26
$ws.Cells.Item(185, 9).Value = 7  # I185: was -1
$ws.Cells.Item(185, 16).Value = 'NA l inferred'  # P185: was None
$ws.Cells.Item(186, 5).Value = 'PAPELON'  # E186: was 'LIOCDEP'
$ws.Cells.Item(186, 8).Value = 22  # H186: was 25
$ws.Cells.Item(186, 9).Value = 7  # I186: was -1
$ws.Cells.Item(186, 16).Value = 'NA l inferred'  # P186: was None
$ws.Cells.Item(187, 8).Value = 24  # H187: was 33
$ws.Cells.Item(187, 10).Value = 'M'  # J187: was 'F'
$ws.Cells.Item(187, 11).ClearContents()  # K187: was 'MEDPF-1'
$ws.Cells.Item(188, 8).Value = 26  # H188: was 27
$ws.Cells.Item(188, 11).Value = 'MEDPF-1'  # K188: was 'MEDPF-2'
$ws.Cells.Item(189, 8).Value = 25  # H189: was 29
$ws.Cells.Item(189, 10).Value = 'F'  # J189: was 'M'
$ws.Cells.Item(189, 11).Value = 'MEDPF-1'  # K189: was None
$ws.Cells.Item(190, 8).Value = 33  # H190: was 30
$ws.Cells.Item(190, 10).Value = 'F'  # J190: was 'M'
$ws.Cells.Item(190, 11).Value = 'MEDPF-1'  # K190: was None
$ws.Cells.Item(191, 8).Value = 27  # H191: was 24
$ws.Cells.Item(191, 10).Value = 'F'  # J191: was 'M'
$ws.Cells.Item(191, 11).Value = 'MEDPF-2'  # K191: was None
$ws.Cells.Item(193, 8).Value = 30  # H193: was 25
$ws.Cells.Item(193, 10).Value = 'M'  # J193: was 'F'
$ws.Cells.Item(193, 11).ClearContents()  # K193: was 'MEDPF-1'
$ws.Cells.Item(194, 10).Value = 'M'  # J194: was 'F'
$ws.Cells.Item(194, 11).ClearContents()  # K194: was 'MEDPF-1'
$ws.Cells.Item(195, 8).Value = 29  # H195: was 22
$ws.Cells.Item(195, 10).Value = 'M'  # J195: was 'F'
$ws.Cells.Item(195, 11).ClearContents()  # K195: was 'MEDPF-1'
$ws.Cells.Item(196, 8).Value = 25  # H196: was 33
$ws.Cells.Item(196, 10).Value = 'F'  # J196: was 'M'
$ws.Cells.Item(196, 11).Value = 'MEDPF-1'  # K196: was None
$ws.Cells.Item(197, 8).Value = 24  # H197: was 28
$ws.Cells.Item(198, 8).Value = 22  # H198: was 26
$ws.Cells.Item(198, 10).Value = 'F'  # J198: was 'M'
$ws.Cells.Item(198, 11).Value = 'MEDPF-1'  # K198: was None
$ws.Cells.Item(199, 8).Value = 33  # H199: was 26
$ws.Cells.Item(200, 8).Value = 28  # H200: was 36
$ws.Cells.Item(200, 11).Value = 'MEDPF-1'  # K200: was 'MEDPF-2'
$ws.Cells.Item(202, 10).Value = 'M'  # J202: was 'F'
$ws.Cells.Item(202, 11).ClearContents()  # K202: was 'MEDPF-1'
$ws.Cells.Item(203, 8).Value = 36  # H203: was 24
$ws.Cells.Item(203, 10).Value = 'F'  # J203: was 'M'
$ws.Cells.Item(203, 11).Value = 'MEDPF-2'  # K203: was None
$ws.Cells.Item(204, 8).Value = 26  # H204: was 30
$ws.Cells.Item(205, 8).Value = 26  # H205: was 27
$ws.Cells.Item(205, 10).Value = 'F'  # J205: was 'M'
$ws.Cells.Item(205, 11).Value = 'MEDPF-1'  # K205: was None
$ws.Cells.Item(206, 8).Value = 24  # H206: was 23
$ws.Cells.Item(206, 10).Value = 'M'  # J206: was 'F'
$ws.Cells.Item(206, 11).ClearContents()  # K206: was 'MEDPF-1'
$ws.Cells.Item(207, 8).Value = 30  # H207: was 22
$ws.Cells.Item(208, 8).Value = 27  # H208: was 26
$ws.Cells.Item(208, 10).Value = 'M'  # J208: was 'F'
$ws.Cells.Item(208, 11).ClearContents()  # K208: was 'MEDPF-1'
$ws.Cells.Item(209, 8).Value = 23  # H209: was 21
$ws.Cells.Item(209, 11).Value = 'MEDPF-1'  # K209: was 'MEDPF-2'
$ws.Cells.Item(210, 8).Value = 22  # H210: was 26
$ws.Cells.Item(210, 10).Value = 'M'  # J210: was 'F'
$ws.Cells.Item(210, 11).ClearContents()  # K210: was 'MEDPF-1'
$ws.Cells.Item(211, 8).Value = 26  # H211: was 27
$ws.Cells.Item(211, 10).Value = 'F'  # J211: was 'M'
$ws.Cells.Item(211, 11).Value = 'MEDPF-1'  # K211: was None
$ws.Cells.Item(212, 8).Value = 21  # H212: was 23
$ws.Cells.Item(212, 10).Value = 'F'  # J212: was 'M'
$ws.Cells.Item(212, 11).Value = 'MEDPF-2'  # K212: was None
$ws.Cells.Item(213, 8).Value = 26  # H213: was 25
$ws.Cells.Item(213, 10).Value = 'F'  # J213: was 'M'
$ws.Cells.Item(213, 11).Value = 'MEDPF-1'  # K213: was None
$ws.Cells.Item(214, 8).Value = 27  # H214: was 31
$ws.Cells.Item(214, 10).Value = 'M'  # J214: was 'F'
$ws.Cells.Item(214, 11).ClearContents()  # K214: was 'MEDPF-1'
$ws.Cells.Item(215, 8).Value = 23  # H215: was 24
$ws.Cells.Item(216, 10).Value = 'M'  # J216: was 'F'
$ws.Cells.Item(216, 11).ClearContents()  # K216: was 'MEDPF-1'
$ws.Cells.Item(217, 8).Value = 31  # H217: was 33
$ws.Cells.Item(217, 10).Value = 'F'  # J217: was 'M'
$ws.Cells.Item(217, 11).Value = 'MEDPF-1'  # K217: was None
$ws.Cells.Item(218, 8).Value = 24  # H218: was 23

# --- Append new rows 219-221 ---
# Row 219: SOLEMON2025 | ITA17 | 38 | 2-RAP | LIOCDEP | 1 | 1 | 25 | -1 | F | MEDPF-1 | 0
$ws.Cells.Item(219, 1).Value = 'SOLEMON2025'
$ws.Cells.Item(219, 2).Value = 'ITA17'
$ws.Cells.Item(218, 3).Copy($ws.Cells.Item(219, 3))
$ws.Cells.Item(219, 4).Value = '2-RAP'
$ws.Cells.Item(219, 5).Value = 'LIOCDEP'
$ws.Cells.Item(219, 6).Value = 1
$ws.Cells.Item(219, 7).Value = 1
$ws.Cells.Item(219, 8).Value = 25
$ws.Cells.Item(219, 9).Value = -1
$ws.Cells.Item(219, 10).Value = 'F'
$ws.Cells.Item(219, 11).Value = 'MEDPF-1'
$ws.Cells.Item(219, 12).Value = 0

# Row 220: SOLEMON2025 | ITA17 | 38 | 2-RAP | LIOCDEP | 1 | 1 | 33 | -1 | M | (no MatStage) | 0
$ws.Cells.Item(220, 1).Value = 'SOLEMON2025'
$ws.Cells.Item(220, 2).Value = 'ITA17'
$ws.Cells.Item(218, 3).Copy($ws.Cells.Item(220, 3))
$ws.Cells.Item(220, 4).Value = '2-RAP'
$ws.Cells.Item(220, 5).Value = 'LIOCDEP'
$ws.Cells.Item(220, 6).Value = 1
$ws.Cells.Item(220, 7).Value = 1
$ws.Cells.Item(220, 8).Value = 33
$ws.Cells.Item(220, 9).Value = -1
$ws.Cells.Item(220, 10).Value = 'M'
$ws.Cells.Item(220, 12).Value = 0

# Row 221: SOLEMON2025 | ITA17 | 38 | 2-RAP | LIOCDEP | 1 | 1 | 23 | -1 | M | (no MatStage) | 0
$ws.Cells.Item(221, 1).Value = 'SOLEMON2025'
$ws.Cells.Item(221, 2).Value = 'ITA17'
$ws.Cells.Item(218, 3).Copy($ws.Cells.Item(221, 3))
$ws.Cells.Item(221, 4).Value = '2-RAP'
$ws.Cells.Item(221, 5).Value = 'LIOCDEP'
$ws.Cells.Item(221, 6).Value = 1
$ws.Cells.Item(221, 7).Value = 1
$ws.Cells.Item(221, 8).Value = 23
$ws.Cells.Item(221, 9).Value = -1
$ws.Cells.Item(221, 10).Value = 'M'
$ws.Cells.Item(221, 12).Value = 0
